# Add a new "2022-Q1" sheet (fund-level detail) positioned right after
# "2021-Q4" and before the existing "总计" summary sheet, then insert a
# matching summary row at the top of "总计" for the new quarter.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet between "2021-Q4" and "总计" ---
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Copy the header/column formatting (bold, centered, bordered - style used
# by the other quarterly detail sheets) from "2021-Q4" onto the new sheet.
$q4.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows (columns B-G stored as text, matching source data)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001759"
$newSheet.Range("C2").Value = "嘉实成长增强灵活配置混合"
$newSheet.Range("D2").Value = "'4.59"
$newSheet.Range("E2").Value = "'90.80"
$newSheet.Range("F2").Value = "'3.82"
$newSheet.Range("G2").Value = "'0.1753"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'005305"
$newSheet.Range("C3").Value = "长信合利混合A"
$newSheet.Range("D3").Value = "'1.99"
$newSheet.Range("E3").Value = "'38.27"
$newSheet.Range("F3").Value = "'1.72"
$newSheet.Range("G3").Value = "'0.0342"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'004608"
$newSheet.Range("C4").Value = "长信乐信灵活配置混合A"
$newSheet.Range("D4").Value = "'0.59"
$newSheet.Range("E4").Value = "'40.07"
$newSheet.Range("F4").Value = "'1.94"
$newSheet.Range("G4").Value = "'0.0114"
$newSheet.Range("H4").Value = 9

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'001744"
$newSheet.Range("C5").Value = "诺安进取回报灵活配置混合"
$newSheet.Range("D5").Value = "'0.04"
$newSheet.Range("E5").Value = "'62.10"
$newSheet.Range("F5").Value = "'4.36"
$newSheet.Range("G5").Value = "'0.0017"
$newSheet.Range("H5").Value = 8

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'004609"
$newSheet.Range("C6").Value = "长信乐信灵活配置混合C"
$newSheet.Range("D6").Value = "'0.04"
$newSheet.Range("E6").Value = "'40.07"
$newSheet.Range("F6").Value = "'1.94"
$newSheet.Range("G6").Value = "'0.0008"
$newSheet.Range("H6").Value = 9

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'005306"
$newSheet.Range("C7").Value = "长信合利混合C"
$newSheet.Range("D7").Value = "'0.00"
$newSheet.Range("E7").Value = "'38.27"
$newSheet.Range("F7").Value = "'1.72"
$newSheet.Range("G7").Value = 0
$newSheet.Range("H7").Value = 10

# --- 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet ---
# (fetched fresh now that the sheet collection/index has shifted because of
# the new sheet inserted above)
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# New row inherits formatting from the row insert; reset the data cells to
# the plain (unstyled) look used by the other data rows, keep column A's
# bordered/centered index style.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.22

# Re-number the index column for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Restore the originally active sheet/tab (unrelated to this edit)
$wb.Worksheets.Item("2021-Q3").Activate()
